$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing mAP (D) and Avg. frame processing time (E) benchmark
# results for the IoU=0.75 and IoU=0.95 blocks (rows 7-16) — the "small,
# 180 (image sample) data benchmark" rows that were previously blank.

# IoU threshold = 0.75
$ws.Range("D7").Value = 0.0592592592592592
$ws.Range("E7").Value = 0.0369614164034525

$ws.Range("D8").Value = 0.504907407407407
$ws.Range("E8").Value = 0.0415858414438035

$ws.Range("D9").Value = 0.476388888888888
$ws.Range("E9").Value = 1.22331076198154

$ws.Range("D10").Value = 0.0379629629629629
$ws.Range("E10").Value = 0.0957023329204983

$ws.Range("D11").Value = 0.0648148148148148
$ws.Range("E11").Value = 13.6222740888595

# IoU threshold = 0.95
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0.036544668674469

$ws.Range("D13").Value = 0.00555555555555555
$ws.Range("E13").Value = 0.0412010563744439

$ws.Range("D14").Value = 0.00555555555555555
$ws.Range("E14").Value = 1.21208487749099

$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0.0962087829907735

$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 13.6158024642202

# Now that column D actually has numbers in it, widen it a bit (matches the
# author widening column D to fit the new mAP values) and move the
# selection down below the table, mirroring where the author clicked next.
$ws.Columns("D").ColumnWidth = 11.2857142857143
$ws.Range("D21").Select() | Out-Null
